$wb = $excel.ActiveWorkbook

# --- PIR sheet (index 2): extend dimension from F170 to F183 ---
$ws = $wb.Worksheets.Item(2)
$rng = $ws.Range("A171:F183")
$rng.NumberFormat = "@"
$ws.Cells.Item(171, 1).Value = "2026-01-28"
$ws.Cells.Item(171, 2).Value = "16:51:18"
$ws.Cells.Item(171, 3).Value = "16:00"
$ws.Cells.Item(171, 4).Value = "Bathroom"
$ws.Cells.Item(171, 5).Value = "No Motion"
$ws.Cells.Item(171, 6).Value = "Inactive"
$ws.Cells.Item(172, 1).Value = "2026-01-28"
$ws.Cells.Item(172, 2).Value = "16:51:19"
$ws.Cells.Item(172, 3).Value = "16:00"
$ws.Cells.Item(172, 4).Value = "Bathroom"
$ws.Cells.Item(172, 5).Value = "No Motion"
$ws.Cells.Item(172, 6).Value = "Inactive"
$ws.Cells.Item(173, 1).Value = "2026-01-28"
$ws.Cells.Item(173, 2).Value = "16:51:23"
$ws.Cells.Item(173, 3).Value = "16:00"
$ws.Cells.Item(173, 4).Value = "Bathroom"
$ws.Cells.Item(173, 5).Value = "No Motion"
$ws.Cells.Item(173, 6).Value = "Inactive"
$ws.Cells.Item(174, 1).Value = "2026-01-28"
$ws.Cells.Item(174, 2).Value = "16:51:28"
$ws.Cells.Item(174, 3).Value = "16:00"
$ws.Cells.Item(174, 4).Value = "Bathroom"
$ws.Cells.Item(174, 5).Value = "No Motion"
$ws.Cells.Item(174, 6).Value = "Inactive"
$ws.Cells.Item(175, 1).Value = "2026-01-28"
$ws.Cells.Item(175, 2).Value = "16:51:33"
$ws.Cells.Item(175, 3).Value = "16:00"
$ws.Cells.Item(175, 4).Value = "Bathroom"
$ws.Cells.Item(175, 5).Value = "No Motion"
$ws.Cells.Item(175, 6).Value = "Inactive"
$ws.Cells.Item(176, 1).Value = "2026-01-28"
$ws.Cells.Item(176, 2).Value = "16:51:38"
$ws.Cells.Item(176, 3).Value = "16:00"
$ws.Cells.Item(176, 4).Value = "Bathroom"
$ws.Cells.Item(176, 5).Value = "No Motion"
$ws.Cells.Item(176, 6).Value = "Inactive"
$ws.Cells.Item(177, 1).Value = "2026-01-28"
$ws.Cells.Item(177, 2).Value = "16:51:44"
$ws.Cells.Item(177, 3).Value = "16:00"
$ws.Cells.Item(177, 4).Value = "Bathroom"
$ws.Cells.Item(177, 5).Value = "No Motion"
$ws.Cells.Item(177, 6).Value = "Inactive"
$ws.Cells.Item(178, 1).Value = "2026-01-28"
$ws.Cells.Item(178, 2).Value = "16:51:48"
$ws.Cells.Item(178, 3).Value = "16:00"
$ws.Cells.Item(178, 4).Value = "Bathroom"
$ws.Cells.Item(178, 5).Value = "No Motion"
$ws.Cells.Item(178, 6).Value = "Inactive"
$ws.Cells.Item(179, 1).Value = "2026-01-28"
$ws.Cells.Item(179, 2).Value = "16:51:53"
$ws.Cells.Item(179, 3).Value = "16:00"
$ws.Cells.Item(179, 4).Value = "Bathroom"
$ws.Cells.Item(179, 5).Value = "No Motion"
$ws.Cells.Item(179, 6).Value = "Inactive"
$ws.Cells.Item(180, 1).Value = "2026-01-28"
$ws.Cells.Item(180, 2).Value = "16:51:58"
$ws.Cells.Item(180, 3).Value = "16:00"
$ws.Cells.Item(180, 4).Value = "Bathroom"
$ws.Cells.Item(180, 5).Value = "No Motion"
$ws.Cells.Item(180, 6).Value = "Inactive"
$ws.Cells.Item(181, 1).Value = "2026-01-28"
$ws.Cells.Item(181, 2).Value = "16:52:03"
$ws.Cells.Item(181, 3).Value = "16:00"
$ws.Cells.Item(181, 4).Value = "Bathroom"
$ws.Cells.Item(181, 5).Value = "No Motion"
$ws.Cells.Item(181, 6).Value = "Inactive"
$ws.Cells.Item(182, 1).Value = "2026-01-28"
$ws.Cells.Item(182, 2).Value = "16:52:08"
$ws.Cells.Item(182, 3).Value = "16:00"
$ws.Cells.Item(182, 4).Value = "Bathroom"
$ws.Cells.Item(182, 5).Value = "No Motion"
$ws.Cells.Item(182, 6).Value = "Inactive"
$ws.Cells.Item(183, 1).Value = "2026-01-28"
$ws.Cells.Item(183, 2).Value = "16:52:13"
$ws.Cells.Item(183, 3).Value = "16:00"
$ws.Cells.Item(183, 4).Value = "Bathroom"
$ws.Cells.Item(183, 5).Value = "No Motion"
$ws.Cells.Item(183, 6).Value = "Inactive"

# --- Humidity sheet (index 3): extend dimension from F167 to F177 ---
$ws = $wb.Worksheets.Item(3)
$rng = $ws.Range("A168:F177")
$rng.NumberFormat = "@"
$ws.Cells.Item(168, 1).Value = "2026-01-28"
$ws.Cells.Item(168, 2).Value = "16:51:16"
$ws.Cells.Item(168, 3).Value = "16:00"
$ws.Cells.Item(168, 4).Value = "Bathroom"
$ws.Cells.Item(168, 5).Value = "86.9%"
$ws.Cells.Item(168, 6).Value = "Active"
$ws.Cells.Item(169, 1).Value = "2026-01-28"
$ws.Cells.Item(169, 2).Value = "16:51:20"
$ws.Cells.Item(169, 3).Value = "16:00"
$ws.Cells.Item(169, 4).Value = "Bathroom"
$ws.Cells.Item(169, 5).Value = "87.9%"
$ws.Cells.Item(169, 6).Value = "Active"
$ws.Cells.Item(170, 1).Value = "2026-01-28"
$ws.Cells.Item(170, 2).Value = "16:51:22"
$ws.Cells.Item(170, 3).Value = "16:00"
$ws.Cells.Item(170, 4).Value = "Bathroom"
$ws.Cells.Item(170, 5).Value = "86.9%"
$ws.Cells.Item(170, 6).Value = "Active"
$ws.Cells.Item(171, 1).Value = "2026-01-28"
$ws.Cells.Item(171, 2).Value = "16:51:26"
$ws.Cells.Item(171, 3).Value = "16:00"
$ws.Cells.Item(171, 4).Value = "Bathroom"
$ws.Cells.Item(171, 5).Value = "87.8%"
$ws.Cells.Item(171, 6).Value = "Active"
$ws.Cells.Item(172, 1).Value = "2026-01-28"
$ws.Cells.Item(172, 2).Value = "16:51:34"
$ws.Cells.Item(172, 3).Value = "16:00"
$ws.Cells.Item(172, 4).Value = "Bathroom"
$ws.Cells.Item(172, 5).Value = "86.9%"
$ws.Cells.Item(172, 6).Value = "Active"
$ws.Cells.Item(173, 1).Value = "2026-01-28"
$ws.Cells.Item(173, 2).Value = "16:51:43"
$ws.Cells.Item(173, 3).Value = "16:00"
$ws.Cells.Item(173, 4).Value = "Bathroom"
$ws.Cells.Item(173, 5).Value = "86.9%"
$ws.Cells.Item(173, 6).Value = "Active"
$ws.Cells.Item(174, 1).Value = "2026-01-28"
$ws.Cells.Item(174, 2).Value = "16:51:47"
$ws.Cells.Item(174, 3).Value = "16:00"
$ws.Cells.Item(174, 4).Value = "Bathroom"
$ws.Cells.Item(174, 5).Value = "87.9%"
$ws.Cells.Item(174, 6).Value = "Active"
$ws.Cells.Item(175, 1).Value = "2026-01-28"
$ws.Cells.Item(175, 2).Value = "16:51:55"
$ws.Cells.Item(175, 3).Value = "16:00"
$ws.Cells.Item(175, 4).Value = "Bathroom"
$ws.Cells.Item(175, 5).Value = "87.9%"
$ws.Cells.Item(175, 6).Value = "Active"
$ws.Cells.Item(176, 1).Value = "2026-01-28"
$ws.Cells.Item(176, 2).Value = "16:52:07"
$ws.Cells.Item(176, 3).Value = "16:00"
$ws.Cells.Item(176, 4).Value = "Bathroom"
$ws.Cells.Item(176, 5).Value = "87.9%"
$ws.Cells.Item(176, 6).Value = "Active"
$ws.Cells.Item(177, 1).Value = "2026-01-28"
$ws.Cells.Item(177, 2).Value = "16:52:11"
$ws.Cells.Item(177, 3).Value = "16:00"
$ws.Cells.Item(177, 4).Value = "Bathroom"
$ws.Cells.Item(177, 5).Value = "87.9%"
$ws.Cells.Item(177, 6).Value = "Active"

# --- Temperature sheet (index 4): extend dimension from F167 to F177 ---
$ws = $wb.Worksheets.Item(4)
$rng = $ws.Range("A168:F177")
$rng.NumberFormat = "@"
$ws.Cells.Item(168, 1).Value = "2026-01-28"
$ws.Cells.Item(168, 2).Value = "16:51:17"
$ws.Cells.Item(168, 3).Value = "16:00"
$ws.Cells.Item(168, 4).Value = "Bathroom"
$ws.Cells.Item(168, 5).Value = "22.8C"
$ws.Cells.Item(168, 6).Value = "Active"
$ws.Cells.Item(169, 1).Value = "2026-01-28"
$ws.Cells.Item(169, 2).Value = "16:51:20"
$ws.Cells.Item(169, 3).Value = "16:00"
$ws.Cells.Item(169, 4).Value = "Bathroom"
$ws.Cells.Item(169, 5).Value = "22.8C"
$ws.Cells.Item(169, 6).Value = "Active"
$ws.Cells.Item(170, 1).Value = "2026-01-28"
$ws.Cells.Item(170, 2).Value = "16:51:23"
$ws.Cells.Item(170, 3).Value = "16:00"
$ws.Cells.Item(170, 4).Value = "Bathroom"
$ws.Cells.Item(170, 5).Value = "22.8C"
$ws.Cells.Item(170, 6).Value = "Active"
$ws.Cells.Item(171, 1).Value = "2026-01-28"
$ws.Cells.Item(171, 2).Value = "16:51:27"
$ws.Cells.Item(171, 3).Value = "16:00"
$ws.Cells.Item(171, 4).Value = "Bathroom"
$ws.Cells.Item(171, 5).Value = "22.8C"
$ws.Cells.Item(171, 6).Value = "Active"
$ws.Cells.Item(172, 1).Value = "2026-01-28"
$ws.Cells.Item(172, 2).Value = "16:51:35"
$ws.Cells.Item(172, 3).Value = "16:00"
$ws.Cells.Item(172, 4).Value = "Bathroom"
$ws.Cells.Item(172, 5).Value = "22.8C"
$ws.Cells.Item(172, 6).Value = "Active"
$ws.Cells.Item(173, 1).Value = "2026-01-28"
$ws.Cells.Item(173, 2).Value = "16:51:43"
$ws.Cells.Item(173, 3).Value = "16:00"
$ws.Cells.Item(173, 4).Value = "Bathroom"
$ws.Cells.Item(173, 5).Value = "22.8C"
$ws.Cells.Item(173, 6).Value = "Active"
$ws.Cells.Item(174, 1).Value = "2026-01-28"
$ws.Cells.Item(174, 2).Value = "16:51:47"
$ws.Cells.Item(174, 3).Value = "16:00"
$ws.Cells.Item(174, 4).Value = "Bathroom"
$ws.Cells.Item(174, 5).Value = "22.8C"
$ws.Cells.Item(174, 6).Value = "Active"
$ws.Cells.Item(175, 1).Value = "2026-01-28"
$ws.Cells.Item(175, 2).Value = "16:51:55"
$ws.Cells.Item(175, 3).Value = "16:00"
$ws.Cells.Item(175, 4).Value = "Bathroom"
$ws.Cells.Item(175, 5).Value = "22.8C"
$ws.Cells.Item(175, 6).Value = "Active"
$ws.Cells.Item(176, 1).Value = "2026-01-28"
$ws.Cells.Item(176, 2).Value = "16:52:07"
$ws.Cells.Item(176, 3).Value = "16:00"
$ws.Cells.Item(176, 4).Value = "Bathroom"
$ws.Cells.Item(176, 5).Value = "22.8C"
$ws.Cells.Item(176, 6).Value = "Active"
$ws.Cells.Item(177, 1).Value = "2026-01-28"
$ws.Cells.Item(177, 2).Value = "16:52:11"
$ws.Cells.Item(177, 3).Value = "16:00"
$ws.Cells.Item(177, 4).Value = "Bathroom"
$ws.Cells.Item(177, 5).Value = "22.8C"
$ws.Cells.Item(177, 6).Value = "Active"

# --- mmWave sheet (index 6): extend dimension from F11 to F25 ---
$ws = $wb.Worksheets.Item(6)
$rng = $ws.Range("A12:F25")
$rng.NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "2026-01-28"
$ws.Cells.Item(12, 2).Value = "16:51:17"
$ws.Cells.Item(12, 3).Value = "16:00"
$ws.Cells.Item(12, 4).Value = "Living Room"
$ws.Cells.Item(12, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(12, 6).Value = "Active"
$ws.Cells.Item(13, 1).Value = "2026-01-28"
$ws.Cells.Item(13, 2).Value = "16:51:18"
$ws.Cells.Item(13, 3).Value = "16:00"
$ws.Cells.Item(13, 4).Value = "Living Room"
$ws.Cells.Item(13, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(13, 6).Value = "Active"
$ws.Cells.Item(14, 1).Value = "2026-01-28"
$ws.Cells.Item(14, 2).Value = "16:51:19"
$ws.Cells.Item(14, 3).Value = "16:00"
$ws.Cells.Item(14, 4).Value = "Living Room"
$ws.Cells.Item(14, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(14, 6).Value = "Active"
$ws.Cells.Item(15, 1).Value = "2026-01-28"
$ws.Cells.Item(15, 2).Value = "16:51:21"
$ws.Cells.Item(15, 3).Value = "16:00"
$ws.Cells.Item(15, 4).Value = "Living Room"
$ws.Cells.Item(15, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(15, 6).Value = "Active"
$ws.Cells.Item(16, 1).Value = "2026-01-28"
$ws.Cells.Item(16, 2).Value = "16:51:24"
$ws.Cells.Item(16, 3).Value = "16:00"
$ws.Cells.Item(16, 4).Value = "Living Room"
$ws.Cells.Item(16, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(16, 6).Value = "Active"
$ws.Cells.Item(17, 1).Value = "2026-01-28"
$ws.Cells.Item(17, 2).Value = "16:51:27"
$ws.Cells.Item(17, 3).Value = "16:00"
$ws.Cells.Item(17, 4).Value = "Living Room"
$ws.Cells.Item(17, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(17, 6).Value = "Active"
$ws.Cells.Item(18, 1).Value = "2026-01-28"
$ws.Cells.Item(18, 2).Value = "16:51:30"
$ws.Cells.Item(18, 3).Value = "16:00"
$ws.Cells.Item(18, 4).Value = "Living Room"
$ws.Cells.Item(18, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(18, 6).Value = "Active"
$ws.Cells.Item(19, 1).Value = "2026-01-28"
$ws.Cells.Item(19, 2).Value = "16:51:34"
$ws.Cells.Item(19, 3).Value = "16:00"
$ws.Cells.Item(19, 4).Value = "Living Room"
$ws.Cells.Item(19, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(19, 6).Value = "Active"
$ws.Cells.Item(20, 1).Value = "2026-01-28"
$ws.Cells.Item(20, 2).Value = "16:51:36"
$ws.Cells.Item(20, 3).Value = "16:00"
$ws.Cells.Item(20, 4).Value = "Living Room"
$ws.Cells.Item(20, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(20, 6).Value = "Active"
$ws.Cells.Item(21, 1).Value = "2026-01-28"
$ws.Cells.Item(21, 2).Value = "16:51:39"
$ws.Cells.Item(21, 3).Value = "16:00"
$ws.Cells.Item(21, 4).Value = "Living Room"
$ws.Cells.Item(21, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(21, 6).Value = "Active"
$ws.Cells.Item(22, 1).Value = "2026-01-28"
$ws.Cells.Item(22, 2).Value = "16:51:42"
$ws.Cells.Item(22, 3).Value = "16:00"
$ws.Cells.Item(22, 4).Value = "Living Room"
$ws.Cells.Item(22, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(22, 6).Value = "Active"
$ws.Cells.Item(23, 1).Value = "2026-01-28"
$ws.Cells.Item(23, 2).Value = "16:51:45"
$ws.Cells.Item(23, 3).Value = "16:00"
$ws.Cells.Item(23, 4).Value = "Living Room"
$ws.Cells.Item(23, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(23, 6).Value = "Active"
$ws.Cells.Item(24, 1).Value = "2026-01-28"
$ws.Cells.Item(24, 2).Value = "16:51:49"
$ws.Cells.Item(24, 3).Value = "16:00"
$ws.Cells.Item(24, 4).Value = "Living Room"
$ws.Cells.Item(24, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(24, 6).Value = "Active"
$ws.Cells.Item(25, 1).Value = "2026-01-28"
$ws.Cells.Item(25, 2).Value = "16:51:51"
$ws.Cells.Item(25, 3).Value = "16:00"
$ws.Cells.Item(25, 4).Value = "Living Room"
$ws.Cells.Item(25, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(25, 6).Value = "Active"
